# todolist.xlsx update - add the rest of the project plan (Website / Consumer / Testing /
# Social Media / Company sections) below the existing Vendor section, fix up the couple of
# Vendor rows whose text changed, and move the selection to the new last cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target text for every shared string the finished sheet references (by the shared-string
# index it ends up at in the diff) - we just assign literal cell values; the engine owns
# shared-string de-duplication / table layout.
$S = @{}
$S[0] = 'Task'
$S[1] = 'Description'
$S[2] = 'Section'
$S[3] = 'Vendor'
$S[4] = 'Complete Vendor Model file'
$S[5] = 'Update all of the functions for the model file for vendors. So all of the data is ready to be pulled.'
$S[6] = 'Track Finish Date'
$S[7] = 'Create Manage Users Page'
$S[8] = 'Create a page to manage additional vendor users for vendors'
$S[9] = 'Create Manage Business Page'
$S[10] = 'Create Manage Business Information Page'
$S[11] = 'Create Manage Promos Page'
$S[12] = 'This page is what initially is loaded for a vendor. It shows listings of all of the businesses they can manage including a manage button. It gets if they are a premium member and their business basic information just for presentation purposes. This also includes editing their business information such as hours, description, and keywords.'
$S[13] = 'This page gets all the promotions and posts they made for consumers feeds. It will allow them to update any text, and image they would like to add. Like a facebook post. It will allow them to delete any old posts they may have made in the past.'
$S[14] = 'Create PPC Campaign Pages'
$S[15] = 'If the user is a premium member alllow them to add credit to their account. As well as create a pay per click campaign with as many keywords as they want. With the ability to delete campaign, make it inactive, add/edit/remove keywords.'
$S[16] = 'Create Reports Page'
$S[17] = 'The reports page will get business review stats, business stats general, ppc stats if applicable, rating stats, mobile app stats, and web stats. All of this is for premium members with the exception of some basic statistical data which counts as business stats general.'
$S[18] = 'This page is once you get into the vendor business. This shows all of the basic options they can do such as edit their business information. As well as add photos and create menu.'
$S[19] = 'Create Manage Reviews Page'
$S[20] = 'This page will go hand and hand with the review system. Allowing users to delete reviews made by people without accounts. Delete 1 review a month of people who are registered but not verified and submit for deletion review for all other registered verified members including non verified incase they went over their limit. This page will also allow them to respond to customers. Customers will be emailed if subscribed to wheres the grub as a response as been made upon their review.'
$S[21] = 'Create a subscription signup page'
$S[22] = 'This page will draw the users to our services as a premium member. It will layout all the qualities of being a premium member then a regular member. It will also have a signup for month to month, 3 month, 6month, and 1 year subscription. As we will promote signing up for 1 year at a time with a recurring subscription billing setup.'
$S[23] = 'Website'
$S[24] = 'Update all listing results'
$S[25] = 'I will need to change the way we get listing results for all parts of the website to effectively show premium and PPC members being ahead of everybody as well as showing up as a competitor.'
$S[26] = 'Update restaurant profiles'
$S[27] = 'I will need to update the restaurant profiles to show anything that’s missing, menu items, review responses, follow button for consumers, photos, etc.'
$S[28] = 'Create where are you page'
$S[29] = 'If I cannot accurately get where somebody is by their zipcode initially I will bounce them to a page that requires them to submit where they are. '
$S[30] = 'Create link tracking'
$S[31] = 'For all vendorstats_types I must update the website so these details are actually tracked in the database for reports and stats'
$S[32] = 'Create services page for restaurants/adding restaurant'
$S[33] = 'Create a page for the public to see what it takes to create a listing as a restaurant to add'
$S[34] = 'Consumer'
$S[35] = 'Update profile page'
$S[36] = 'Make sure users can upload photos of themselves as profile pictures for their avatar, update their bio, remove followed restaurants, add followed restaurants. '
$S[37] = 'Update feed page'
$S[38] = 'Create the feed page for information of coupons, deals, of people the consumer followed in chronological order.'
$S[39] = 'Reviews Check'
$S[40] = 'Make sure consumers can leave reviews, get activated accounts, and that anonymous people can leave reviews on restaurant page.'
$S[41] = 'Testing'
$S[42] = 'Test the website as a whole'
$S[43] = 'Test each page, consumer, vendor, website. Test all functionalities. Go over everything and anything. Double check for bugs. '
$S[44] = 'Social Media'
$S[45] = 'Update all social medias'
$S[46] = 'Submit announcements that we will be launching soon. Specify date, update all additional information needed on social platforms.'
$S[47] = 'Company'
$S[48] = 'Soft Launch'
$S[49] = 'We will be launching our website into the public atmosphere of the interwebs. This is a soft launch so we will not be doing extreme mass marketing to promote the site until we get some initial feedback. Marketing comes with dates listed below in time schedule.'

# --- Rows 5-9: existing Vendor rows whose Task/Description text changed ---
# --- Rows 10-28: brand-new rows (rows 12, 18, 22, 25, 26, 27 stay blank, matching the diff) ---
$rowCols = @{
  5 = @{ A = 3; B = 9; C = 12 }
  6 = @{ A = 3; B = 10; C = 18 }
  7 = @{ A = 3; B = 11; C = 13 }
  8 = @{ A = 3; B = 14; C = 15 }
  9 = @{ A = 3; B = 16; C = 17 }
  10 = @{ A = 3; B = 19; C = 20 }
  11 = @{ A = 3; B = 21; C = 22 }
  13 = @{ A = 23; B = 24; C = 25 }
  14 = @{ A = 23; B = 26; C = 27 }
  15 = @{ A = 23; B = 28; C = 29 }
  16 = @{ A = 23; B = 30; C = 31 }
  17 = @{ A = 23; B = 32; C = 33 }
  19 = @{ A = 34; B = 35; C = 36 }
  20 = @{ A = 34; B = 37; C = 38 }
  21 = @{ A = 34; B = 39; C = 40 }
  23 = @{ A = 41; B = 42; C = 43 }
  24 = @{ A = 44; B = 45; C = 46 }
  28 = @{ A = 47; B = 48; C = 49 }
}

$dateSerials = @{
  5 = 42736
  6 = 42738
  7 = 42739
  8 = 42743
  9 = 42748
  10 = 42752
  11 = 42756
  13 = 42757
  14 = 42759
  15 = 42760
  16 = 42761
  17 = 42762
  19 = 42763
  20 = 42764
  21 = 42765
  23 = 42766
  24 = 42766
  28 = 42767
}

foreach ($r in ($rowCols.Keys | Sort-Object)) {
  $cols = $rowCols[$r]
  $ws.Range("A$r").Value = $S[$cols.A]
  $ws.Range("B$r").Value = $S[$cols.B]
  $ws.Range("C$r").Value = $S[$cols.C]

  # Write the finish-date as a bare serial so no ad-hoc number format gets minted, then
  # copy the date format from an untouched date cell (D3) so it lands on the same style
  # index (s="1") the rest of the Track Finish Date column already uses.
  $ws.Range("D$r").Value = $dateSerials[$r]
  $ws.Range("D3").Copy()
  $ws.Range("D$r").PasteSpecial(-4122)
}

# C5 keeps its original "vertical top" style (s="2") - only its text target changed, which
# the loop above already wrote; nothing else to do for formatting there.

# New used range is A1:D28, and Excel leaves the selection on the next empty row, D29.
$ws.Range("D29").Select()

